$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 204
$ws.Range("F3").Value = 1084
$ws.Range("F4").Value = 1131
$ws.Range("F7").Value = 730
$ws.Range("F8").Value = 13174
$ws.Range("F9").Value = 13174
$ws.Range("F10").Value = 2274
$ws.Range("F12").Value = 304
$ws.Range("F13").Value = 54096
$ws.Range("F14").Value = 1299
$ws.Range("F15").Value = 318
$ws.Range("F16").Value = 305
$ws.Range("F17").Value = 857
$ws.Range("F18").Value = 708
$ws.Range("F19").Value = 358
$ws.Range("F20").Value = 2977
$ws.Range("F21").Value = 863
$ws.Range("F22").Value = 5140
$ws.Range("F23").Value = 1252
$ws.Range("F24").Value = 938
$ws.Range("F26").Value = 35
$ws.Range("F28").Value = 375
$ws.Range("F29").Value = 1201
$ws.Range("F32").Value = 147
$ws.Range("F33").Value = 333
$ws.Range("F34").Value = 36
$ws.Range("F35").Value = 22
$ws.Range("F36").Value = 60
$ws.Range("F37").Value = 46
$ws.Range("F38").Value = 4731
$ws.Range("F39").Value = 34
$ws.Range("F40").Value = 4750
$ws.Range("F41").Value = 8713
$ws.Range("F42").Value = 108
$ws.Range("F44").Value = 119
$ws.Range("F45").Value = 205
$ws.Range("F47").Value = 100
$ws.Range("F49").Value = 4164
$ws.Range("F50").Value = 181

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 90
$ws.Range("F12").Value = 1116

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 781
$ws.Range("F5").Value = 32

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 781
$ws.Range("F4").Value = 204
$ws.Range("F5").Value = 1084
$ws.Range("F6").Value = 1131
$ws.Range("F8").Value = 730
$ws.Range("F9").Value = 13174
$ws.Range("F10").Value = 2274
$ws.Range("F11").Value = 1299
$ws.Range("F12").Value = 305
$ws.Range("F13").Value = 857
$ws.Range("F14").Value = 708
$ws.Range("F15").Value = 2977
$ws.Range("F16").Value = 863
$ws.Range("F17").Value = 90
$ws.Range("F18").Value = 1252
$ws.Range("F19").Value = 32
$ws.Range("F20").Value = 938
$ws.Range("F23").Value = 35
$ws.Range("F28").Value = 147
$ws.Range("F29").Value = 333
$ws.Range("F30").Value = 36
$ws.Range("F31").Value = 46
$ws.Range("F32").Value = 4731
$ws.Range("F33").Value = 34
$ws.Range("F34").Value = 4750
$ws.Range("F36").Value = 108
$ws.Range("F38").Value = 119
$ws.Range("F39").Value = 205
$ws.Range("F43").Value = 100
$ws.Range("F45").Value = 4164
